# Domino Single USB Rev. C: changed Logos
# Rev. B -> Rev. C sheet rename, print-area bookkeeping, selection & column width refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldName = $ws.Name
$newName = "Domino Single USB Rev. C"

# --- Rename the worksheet (cascades through most of the _xlnm.Print_Area_* names) ---
$ws.Name = $newName

# --- The bare "_xlnm.Print_Area" entry doesn't get updated by the rename, and Excel's
#     history for this sheet carries a duplicate of it (two identical defined names
#     both called Print_Area, scoped to this sheet). Fix the stale one in place ... ---
$fixed = $false
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Print_Area" -and -not $fixed) {
        $n.RefersTo = "='" + $newName + "'!`$A`$1:`$I`$10"
        $fixed = $true
    }
}

# --- ... then restore the second (duplicate) Print_Area entry and append one more
#     generation of the _0 history chain, matching the workbook's print-area log. ---
$ws.Names.Add("_xlnm.Print_Area", "='" + $newName + "'!`$A`$1:`$I`$10")
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0", "='" + $newName + "'!`$A`$1:`$I`$1")

# --- Selection narrows from A2:H10 down to just A2 ---
$ws.Range("A2").Select()

# --- Minor column width refresh (font metrics re-measured) ---
$ws.Range("A1").ColumnWidth = 4
$ws.Range("B1").ColumnWidth = 4
$ws.Range("C1").ColumnWidth = 25.333333333333336
$ws.Range("D1").ColumnWidth = 27.166666666666668
$ws.Range("E1").ColumnWidth = 30.333333333333336
$ws.Range("F1").ColumnWidth = 25.333333333333336
$ws.Range("G1").ColumnWidth = 39.33333333333333
$ws.Range("H1").ColumnWidth = 61.166666666666664
$ws.Range("I1").ColumnWidth = 23.666666666666668
